$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells stay text (avoid Excel auto-numeric coercion)

# Row 2: Bitcoin -> Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.457.94"
$ws.Range("E2").Value = "  -3.86%  "

# Row 3: Ethereum -> Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.565.86"
$ws.Range("E3").Value = "  -4.45%  "

# Row 4: TetherUSD -> TetherUSD
$ws.Range("E4").Value = "  +0.10%  "

# Row 5: BNB -> BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.78"
$ws.Range("E5").Value = "  -4.31%  "

# Row 6: Solana -> Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.95"
$ws.Range("E6").Value = "  -0.44%  "

# Row 7: LidoStakedEther -> LidoStakedEther
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.561.45"
$ws.Range("E7").Value = "  -4.49%  "

# Row 8: XRP -> XRP
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.614"
$ws.Range("E8").Value = "  -4.20%  "

# Row 9: USDC -> USDC
$ws.Range("E9").Value = "  +0.22%  "

# Row 10: Cardano -> Cardano
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.671"
$ws.Range("E10").Value = "  -6.95%  "

# Row 11: Dogecoin -> Dogecoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.146"
$ws.Range("E11").Value = "  -9.98%  "

# Row 12: Avalanche -> Avalanche
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.28"
$ws.Range("E12").Value = "  -6.13%  "

# Row 13: ShibaInu -> ShibaInu
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000260"
$ws.Range("E13").Value = "  -11.46%  "

# Row 14: Polkadot -> Polkadot
$ws.Range("E14").Value = "  -7.76%  "

# Row 15: WrappedliquidstakedEther2.0 -> WrappedliquidstakedEther2.0
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.140.39"
$ws.Range("E15").Value = "  -4.17%  "

# Row 16: WrappedEther -> WrappedEther
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.566.86"
$ws.Range("E16").Value = "  -4.50%  "

# Row 17: TRON -> TRON
$ws.Range("E17").Value = "  -0.91%  "

# Row 18: Chainlink -> Chainlink
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.35"
$ws.Range("E18").Value = "  -5.36%  "

# Row 19: Uniswap -> Uniswap
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.24"
$ws.Range("E19").Value = "  -6.27%  "

# Row 20: WrappedBTC -> WrappedBTC
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "66.382.21"
$ws.Range("E20").Value = "  -3.72%  "

# Row 21: Polygon -> Polygon
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.05"
$ws.Range("E21").Value = "  -7.63%  "

# Row 22: BitcoinCash -> BitcoinCash
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "395.66"
$ws.Range("E22").Value = "  -4.58%  "

# Row 23: PancakeSwap -> PancakeSwap
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.36"
$ws.Range("E23").Value = "  -5.75%  "

# Row 24: Litecoin -> Litecoin
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "86.09"
$ws.Range("E24").Value = "  -3.68%  "

# Row 25: RenderToken -> RenderToken
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.32"
$ws.Range("E25").Value = "  +2.88%  "

# Row 26: ImmutableX -> ImmutableX
$ws.Range("E26").Value = "  -4.91%  "

# Row 27: InternetComputer(DFINITY) -> InternetComputer(DFINITY)
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.48"
$ws.Range("E27").Value = "  -3.43%  "

# Row 28: LEO -> LEO
$ws.Range("E28").Value = "  -0.29%  "

# Row 29: Toncoin -> Toncoin
$ws.Range("E29").Value = "  -6.41%  "

# Row 30: Filecoin -> Filecoin
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.95"
$ws.Range("E30").Value = "  -7.56%  "

# Row 31: EthereumClassic -> EthereumClassic
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "31.11"
$ws.Range("E31").Value = "  -6.82%  "

# Row 32: NEARProtocol -> NEARProtocol
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.15"
$ws.Range("E32").Value = "  -2.55%  "

# Row 33: Cosmos -> Cosmos
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.18"
$ws.Range("E33").Value = "  -4.72%  "

# Row 34: Bittensor -> Bittensor
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "623.05"
$ws.Range("E34").Value = "  -0.80%  "

# Row 35: Hedera -> Hedera
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.114"
$ws.Range("E35").Value = "  -7.70%  "

# Row 36: OKB -> OKB
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "63.42"
$ws.Range("E36").Value = "  -4.06%  "

# Row 37: InjectiveProtocol -> InjectiveProtocol
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "41.52"
$ws.Range("E37").Value = "  -7.40%  "

# Row 38: Dai -> TheGraph
$ws.Range("B38").Value = "TheGraph"
$ws.Range("C38").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.404"
$ws.Range("E38").Value = "  -1.99%  "

# Row 39: TheGraph -> Dai
$ws.Range("B39").Value = "Dai"
$ws.Range("C39").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  +0.10%  "

# Row 40: PEPE -> PEPE
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0764"
$ws.Range("E40").Value = "  -10.30%  "

# Row 41: Kaspa -> Kaspa
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.133"
$ws.Range("E41").Value = "  -5.43%  "

# Row 42: FirstDigitalUSD -> FirstDigitalUSD
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  +0.01%  "

# Row 43: Maker -> Maker
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.017.31"
$ws.Range("E43").Value = "  +6.26%  "

# Row 44: ThetaToken -> ThetaToken
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.84"
$ws.Range("E44").Value = "  -7.59%  "

# Row 45: Fetch.AI -> Fetch.AI
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.53"
$ws.Range("E45").Value = "  -4.86%  "

# Row 46: VeChain -> VeChain
$ws.Range("E46").Value = "  -8.04%  "

# Row 47: Stellar -> ApeXProtocol
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.14"
$ws.Range("E47").Value = "  +1.83%  "

# Row 48: ApeXProtocol -> Stellar
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.131"
$ws.Range("E48").Value = "  -7.50%  "

# Row 49: THORChain -> THORChain
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.58"
$ws.Range("E49").Value = "  -6.72%  "

# Row 50: Monero -> Monero
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "137.44"
$ws.Range("E50").Value = "  -3.09%  "

# Row 51: Stacks -> Stacks
$ws.Range("E51").Value = "  -1.25%  "
